# chore: update Sheets via scheduled runner
# Refreshes the cached market-board figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) on each Leve-profit sheet with freshly
# scraped values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 770940.9399999999
$ws.Range("J19").Value = 2310.875
$ws.Range("L19").Value = 2310.875
$ws.Range("N19").Value = -2660.875
$ws.Range("H43").Value = 1525
$ws.Range("I43").Value = 1193.6666
$ws.Range("J43").Value = 1635.4445
$ws.Range("K43").Value = 1193.6666
$ws.Range("L43").Value = 1635.4445
$ws.Range("M43").Value = -1124.6666
$ws.Range("N43").Value = -1773.4445
$ws.Range("H55").Value = 296.53333
$ws.Range("J55").Value = 338.75
$ws.Range("L55").Value = 338.75
$ws.Range("N55").Value = -766.75
$ws.Range("H100").Value = 976.38464
$ws.Range("I100").Value = 766.0833
$ws.Range("K100").Value = 766.0833
$ws.Range("M100").Value = -225.0833
$ws.Range("H129").Value = 899.2192
$ws.Range("I129").Value = 1219.8
$ws.Range("J129").Value = 875.64703
$ws.Range("K129").Value = 3659.4
$ws.Range("L129").Value = 2626.94109
$ws.Range("M129").Value = 1340.6
$ws.Range("N129").Value = -12626.94109
$ws.Range("H137").Value = 1559.6154
$ws.Range("I137").Value = 1427.6
$ws.Range("J137").Value = 1999.6666
$ws.Range("K137").Value = 4282.799999999999
$ws.Range("L137").Value = 5998.9998
$ws.Range("M137").Value = -1732.799999999999
$ws.Range("N137").Value = -11098.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3347.426
$ws.Range("I32").Value = 2248.6042
$ws.Range("K32").Value = 2248.6042
$ws.Range("M32").Value = -1961.6042
$ws.Range("H61").Value = 3599.258
$ws.Range("J61").Value = 10008.5
$ws.Range("L61").Value = 10008.5
$ws.Range("N61").Value = -10432.5
$ws.Range("H74").Value = 1019.9429
$ws.Range("I74").Value = 832.9259
$ws.Range("J74").Value = 1651.125
$ws.Range("K74").Value = 832.9259
$ws.Range("L74").Value = 1651.125
$ws.Range("M74").Value = 41.07410000000004
$ws.Range("N74").Value = -3399.125
$ws.Range("H77").Value = 1019.9429
$ws.Range("I77").Value = 832.9259
$ws.Range("J77").Value = 1651.125
$ws.Range("K77").Value = 4164.6295
$ws.Range("L77").Value = 8255.625
$ws.Range("M77").Value = 203.3705
$ws.Range("N77").Value = -16991.625
$ws.Range("H88").Value = 3963.182
$ws.Range("I88").Value = 2466.3333
$ws.Range("J88").Value = 4524.5
$ws.Range("K88").Value = 2466.3333
$ws.Range("L88").Value = 4524.5
$ws.Range("M88").Value = -2060.3333
$ws.Range("N88").Value = -5336.5
$ws.Range("H91").Value = 3963.182
$ws.Range("I91").Value = 2466.3333
$ws.Range("J91").Value = 4524.5
$ws.Range("K91").Value = 2466.3333
$ws.Range("L91").Value = 4524.5
$ws.Range("M91").Value = -1062.3333
$ws.Range("N91").Value = -7332.5
$ws.Range("H132").Value = 1449.5366
$ws.Range("J132").Value = 2025.579
$ws.Range("L132").Value = 6076.737
$ws.Range("N132").Value = -11136.737
$ws.Range("H136").Value = 3599.258
$ws.Range("J136").Value = 10008.5
$ws.Range("L136").Value = 30025.5
$ws.Range("N136").Value = -35125.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2356.6736
$ws.Range("I31").Value = 1523.0286
$ws.Range("K31").Value = 1523.0286
$ws.Range("M31").Value = -1228.0286
$ws.Range("H34").Value = 2356.6736
$ws.Range("I34").Value = 1523.0286
$ws.Range("K34").Value = 1523.0286
$ws.Range("M34").Value = -1321.0286
$ws.Range("H58").Value = 1450793.4
$ws.Range("I58").Value = 3345607.5
$ws.Range("J58").Value = 1817.9412
$ws.Range("K58").Value = 3345607.5
$ws.Range("L58").Value = 1817.9412
$ws.Range("M58").Value = -3345404.5
$ws.Range("N58").Value = -2223.9412
$ws.Range("H132").Value = 2646.6
$ws.Range("J132").Value = 3554.6365
$ws.Range("L132").Value = 10663.9095
$ws.Range("N132").Value = -15723.9095
$ws.Range("H134").Value = 999.7692
$ws.Range("J134").Value = 999.75
$ws.Range("L134").Value = 2999.25
$ws.Range("N134").Value = -8069.25
$ws.Range("H136").Value = 1450793.4
$ws.Range("I136").Value = 3345607.5
$ws.Range("J136").Value = 1817.9412
$ws.Range("K136").Value = 10036822.5
$ws.Range("L136").Value = 5453.8236
$ws.Range("M136").Value = -10034272.5
$ws.Range("N136").Value = -10553.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 348
$ws.Range("J23").Value = 348
$ws.Range("L23").Value = 1044
$ws.Range("N23").Value = -1514
$ws.Range("H44").Value = 3000
$ws.Range("I44").Value = 3000
$ws.Range("K44").Value = 9000
$ws.Range("M44").Value = -8602
$ws.Range("H63").Value = 900
$ws.Range("I63").Value = 900
$ws.Range("K63").Value = 2700
$ws.Range("M63").Value = -1951
$ws.Range("H66").Value = 900
$ws.Range("I66").Value = 900
$ws.Range("K66").Value = 8100
$ws.Range("M66").Value = -4356
$ws.Range("H68").Value = 680
$ws.Range("I68").Value = 680
$ws.Range("K68").Value = 2040
$ws.Range("M68").Value = -1229
$ws.Range("H71").Value = 680
$ws.Range("I71").Value = 680
$ws.Range("K71").Value = 6120
$ws.Range("M71").Value = -2064

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1716313.6
$ws.Range("I126").Value = 2317492
$ws.Range("K126").Value = 6952476
$ws.Range("M126").Value = -6950006

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3541.2856
$ws.Range("I16").Value = 4149
$ws.Range("J16").Value = 2731
$ws.Range("K16").Value = 4149
$ws.Range("L16").Value = 2731
$ws.Range("M16").Value = -3979
$ws.Range("N16").Value = -3071
$ws.Range("H122").Value = 5709.909
$ws.Range("I122").Value = 2800.6667
$ws.Range("J122").Value = 9201
$ws.Range("K122").Value = 8402.000100000001
$ws.Range("L122").Value = 27603
$ws.Range("M122").Value = -5952.000100000001
$ws.Range("N122").Value = -32503
$ws.Range("H136").Value = 3328.5334
$ws.Range("I136").Value = 1616
$ws.Range("J136").Value = 5285.7144
$ws.Range("K136").Value = 4848
$ws.Range("L136").Value = 15857.1432
$ws.Range("M136").Value = -2298
$ws.Range("N136").Value = -20957.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 16174.75
$ws.Range("J39").Value = 19899.666
$ws.Range("L39").Value = 19899.666
$ws.Range("N39").Value = -20725.666
$ws.Range("H107").Value = 654.2727
$ws.Range("I107").Value = 489.07144
$ws.Range("J107").Value = 943.375
$ws.Range("K107").Value = 1467.21432
$ws.Range("L107").Value = 2830.125
$ws.Range("M107").Value = 452.78568
$ws.Range("N107").Value = -6670.125
$ws.Range("H132").Value = 1888.9259
$ws.Range("I132").Value = 1478.3914
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 4435.174199999999
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -1905.174199999999
$ws.Range("N132").Value = -17808.5
$ws.Range("H136").Value = 16837342
$ws.Range("I136").Value = 24156646
$ws.Range("K136").Value = 72469938
$ws.Range("M136").Value = -72467388
